$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("C3").Value2
